$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 3380.3333
$ws.Cells.Item(6, 9).Value = 71
$ws.Cells.Item(6, 10).Value = 9999
$ws.Cells.Item(6, 11).Value = 213
$ws.Cells.Item(6, 12).Value = 29997
$ws.Cells.Item(6, 13).Value = -101
$ws.Cells.Item(6, 14).Value = -30221
$ws.Cells.Item(11, 8).Value = 3000
$ws.Cells.Item(11, 9).Value = 3000
$ws.Cells.Item(11, 11).Value = 3000
$ws.Cells.Item(11, 13).Value = -2860
$ws.Cells.Item(15, 8).Value = 513.6842
$ws.Cells.Item(15, 9).Value = 513.6842
$ws.Cells.Item(15, 11).Value = 1541.0526
$ws.Cells.Item(15, 13).Value = -1372.0526
$ws.Cells.Item(20, 8).Value = 836.25
$ws.Cells.Item(20, 9).Value = 281.66666
$ws.Cells.Item(20, 10).Value = 2500
$ws.Cells.Item(20, 11).Value = 281.66666
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 13).Value = -51.66665999999998
$ws.Cells.Item(20, 14).Value = -2960
$ws.Cells.Item(33, 8).Value = 641.05554
$ws.Cells.Item(33, 9).Value = 178.58333
$ws.Cells.Item(33, 10).Value = 1566
$ws.Cells.Item(33, 11).Value = 178.58333
$ws.Cells.Item(33, 12).Value = 1566
$ws.Cells.Item(33, 13).Value = 50.41667000000001
$ws.Cells.Item(33, 14).Value = -2024
$ws.Cells.Item(35, 8).Value = 836.25
$ws.Cells.Item(35, 9).Value = 281.66666
$ws.Cells.Item(35, 10).Value = 2500
$ws.Cells.Item(35, 11).Value = 281.66666
$ws.Cells.Item(35, 12).Value = 2500
$ws.Cells.Item(35, 13).Value = 97.33334000000002
$ws.Cells.Item(35, 14).Value = -3258
$ws.Cells.Item(62, 8).Value = 6287.5
$ws.Cells.Item(62, 9).Value = 5781.6
$ws.Cells.Item(62, 10).Value = 7130.6665
$ws.Cells.Item(62, 11).Value = 5781.6
$ws.Cells.Item(62, 12).Value = 7130.6665
$ws.Cells.Item(62, 13).Value = -5157.6
$ws.Cells.Item(62, 14).Value = -8378.6665
$ws.Cells.Item(65, 8).Value = 6287.5
$ws.Cells.Item(65, 9).Value = 5781.6
$ws.Cells.Item(65, 10).Value = 7130.6665
$ws.Cells.Item(65, 11).Value = 28908
$ws.Cells.Item(65, 12).Value = 35653.3325
$ws.Cells.Item(65, 13).Value = -25788
$ws.Cells.Item(65, 14).Value = -41893.3325
$ws.Cells.Item(69, 8).Value = 5100
$ws.Cells.Item(69, 9).Value = 3750
$ws.Cells.Item(69, 11).Value = 11250
$ws.Cells.Item(69, 13).Value = -10376
$ws.Cells.Item(72, 8).Value = 5100
$ws.Cells.Item(72, 9).Value = 3750
$ws.Cells.Item(72, 11).Value = 33750
$ws.Cells.Item(72, 13).Value = -29382
$ws.Cells.Item(80, 8).Value = 1936.8334
$ws.Cells.Item(80, 9).Value = 861.5
$ws.Cells.Item(80, 11).Value = 2584.5
$ws.Cells.Item(80, 13).Value = -1586.5
$ws.Cells.Item(83, 8).Value = 1936.8334
$ws.Cells.Item(83, 9).Value = 861.5
$ws.Cells.Item(83, 11).Value = 7753.5
$ws.Cells.Item(83, 13).Value = -2761.5
$ws.Cells.Item(135, 8).Value = 4000
$ws.Cells.Item(135, 9).Value = 4000
$ws.Cells.Item(135, 11).Value = 36000
$ws.Cells.Item(135, 13).Value = -33465
$ws.Cells.Item(138, 8).Value = 5435.591
$ws.Cells.Item(138, 10).Value = 5766.0835
$ws.Cells.Item(138, 12).Value = 17298.2505
$ws.Cells.Item(138, 14).Value = -27578.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 66609.164
$ws.Cells.Item(24, 10).Value = 66609.164
$ws.Cells.Item(24, 12).Value = 66609.164
$ws.Cells.Item(24, 14).Value = -67357.164
$ws.Cells.Item(32, 8).Value = 3140.139
$ws.Cells.Item(32, 9).Value = 2017.7059
$ws.Cells.Item(32, 11).Value = 2017.7059
$ws.Cells.Item(32, 13).Value = -1730.7059
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 12).Value = ""
$ws.Cells.Item(56, 14).Value = 0
$ws.Cells.Item(74, 8).Value = 1506.0769
$ws.Cells.Item(74, 9).Value = 1455.4
$ws.Cells.Item(74, 11).Value = 1455.4
$ws.Cells.Item(74, 13).Value = -581.4000000000001
$ws.Cells.Item(77, 8).Value = 1506.0769
$ws.Cells.Item(77, 9).Value = 1455.4
$ws.Cells.Item(77, 11).Value = 7277
$ws.Cells.Item(77, 13).Value = -2909
$ws.Cells.Item(100, 8).Value = 66609.164
$ws.Cells.Item(100, 10).Value = 66609.164
$ws.Cells.Item(100, 12).Value = 66609.164
$ws.Cells.Item(100, 14).Value = -68773.164

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2994.125
$ws.Cells.Item(31, 9).Value = 2994.125
$ws.Cells.Item(31, 11).Value = 2994.125
$ws.Cells.Item(31, 13).Value = -2699.125
$ws.Cells.Item(34, 8).Value = 2994.125
$ws.Cells.Item(34, 9).Value = 2994.125
$ws.Cells.Item(34, 11).Value = 2994.125
$ws.Cells.Item(34, 13).Value = -2792.125
$ws.Cells.Item(69, 8).Value = 5661.4
$ws.Cells.Item(69, 9).Value = 5661.4
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 5661.4
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = ""
$ws.Cells.Item(69, 14).Value = -4912.4
$ws.Cells.Item(72, 8).Value = 5661.4
$ws.Cells.Item(72, 9).Value = 5661.4
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 16984.2
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).Value = ""
$ws.Cells.Item(72, 14).Value = -13240.2
$ws.Cells.Item(134, 8).Value = 2637.2354
$ws.Cells.Item(134, 9).Value = 2279.9285
$ws.Cells.Item(134, 11).Value = 6839.7855
$ws.Cells.Item(134, 13).Value = -4304.7855

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 999
$ws.Cells.Item(4, 9).Value = 999
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 2997
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = ""
$ws.Cells.Item(4, 14).Value = -2885
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).Value = ""
$ws.Cells.Item(76, 8).Value = 5525
$ws.Cells.Item(76, 9).Value = 5525
$ws.Cells.Item(76, 11).Value = 16575
$ws.Cells.Item(76, 13).Value = -16192
$ws.Cells.Item(79, 8).Value = 5525
$ws.Cells.Item(79, 9).Value = 5525
$ws.Cells.Item(79, 11).Value = 16575
$ws.Cells.Item(79, 13).Value = -15249
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 13).Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(53, 8).Value = 15000
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 15000
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = ""
$ws.Cells.Item(53, 13).Value = 15000
$ws.Cells.Item(53, 14).Value = -16262
$ws.Cells.Item(102, 8).Value = 1131
$ws.Cells.Item(102, 9).Value = 1131
$ws.Cells.Item(102, 11).Value = 1131
$ws.Cells.Item(102, 13).Value = 491
$ws.Cells.Item(122, 8).Value = 500999.5
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 14).Value = ""
$ws.Cells.Item(132, 8).Value = 4724.5
$ws.Cells.Item(132, 9).Value = 5500
$ws.Cells.Item(132, 11).Value = 16500
$ws.Cells.Item(132, 13).Value = -13970

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 956.4
$ws.Cells.Item(16, 9).Value = 1073.5
$ws.Cells.Item(16, 11).Value = 1073.5
$ws.Cells.Item(16, 13).Value = -903.5
$ws.Cells.Item(46, 8).Value = 2950
$ws.Cells.Item(46, 10).Value = 2950
$ws.Cells.Item(46, 12).Value = 2950
$ws.Cells.Item(46, 14).Value = -3326
$ws.Cells.Item(61, 8).Value = 5000
$ws.Cells.Item(61, 9).Value = 5000
$ws.Cells.Item(61, 11).Value = 5000
$ws.Cells.Item(61, 13).Value = -4798
$ws.Cells.Item(68, 8).Value = 5079.8
$ws.Cells.Item(68, 10).Value = 5079.8
$ws.Cells.Item(68, 12).Value = 5079.8
$ws.Cells.Item(68, 14).Value = -6577.8
$ws.Cells.Item(71, 8).Value = 5079.8
$ws.Cells.Item(71, 10).Value = 5079.8
$ws.Cells.Item(71, 12).Value = 25399
$ws.Cells.Item(71, 14).Value = -32887
$ws.Cells.Item(82, 8).Value = 4222.222
$ws.Cells.Item(82, 9).Value = 2714.2856
$ws.Cells.Item(82, 10).Value = 9500
$ws.Cells.Item(82, 11).Value = 2714.2856
$ws.Cells.Item(82, 12).Value = 9500
$ws.Cells.Item(82, 13).Value = -2353.2856
$ws.Cells.Item(82, 14).Value = -10222
$ws.Cells.Item(85, 8).Value = 4222.222
$ws.Cells.Item(85, 9).Value = 2714.2856
$ws.Cells.Item(85, 10).Value = 9500
$ws.Cells.Item(85, 11).Value = 2714.2856
$ws.Cells.Item(85, 12).Value = 9500
$ws.Cells.Item(85, 13).Value = -1466.2856
$ws.Cells.Item(85, 14).Value = -11996
$ws.Cells.Item(113, 8).Value = 5000
$ws.Cells.Item(113, 9).Value = 5000
$ws.Cells.Item(113, 11).Value = 5000
$ws.Cells.Item(113, 13).Value = -2830

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2000
$ws.Cells.Item(126, 9).Value = 2000
$ws.Cells.Item(126, 11).Value = 6000
$ws.Cells.Item(126, 13).Value = -3530
